# Apply the "create API directory / create checking API" edit:
#  - fix the double-space typo in A10 ("get  videos by category" -> "get videos by category")
#  - replace the old "get/send *bytes" block (M1,M3,M5,M7 / N1,N3,N5,N7) with a new
#    "check" API block (M1:N6) highlighted in red, leaving M7:N7 blank
#  - widen the new M/N columns
#  - move the active selection to M12, matching the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix typo in column A ---
$ws.Range("A10").Value = "get videos by category"

# --- new column widths for M (13) and N (14) ---
# (COM ColumnWidth is stored with a +5/6 character offset and rounded to the
#  nearest 1/6 of a character by this runtime, so feed it the inverse value)
$ws.Columns.Item(13).ColumnWidth = 16.721354166666668
$ws.Columns.Item(14).ColumnWidth = 3.2760416666666665

# --- rebuild the M/N "check" table ---
$ws.Range("M1").Value = "check"

$ws.Range("M2").Value = "check email"
$ws.Range("N2").Value = 51

$ws.Range("M3").Value = "check username"
$ws.Range("N3").Value = 52

$ws.Range("M4").Value = "check user exists"
$ws.Range("N4").Value = 53

$ws.Range("M5").Value = "check video saving"
$ws.Range("N5").Value = 54

$ws.Range("M6").Value = " check playlist saving"
$ws.Range("N6").Value = 55

# Rows 2-6 (the new "check" entries) get the new red fill.
$ws.Range("M2:N6").Interior.Color = 255

# Row 7 used to hold "send image bytes" / 8 - it's now blank with no fill.
$ws.Range("M7:N7").Value = $null
$ws.Range("M7:N7").Style = "Normal"

# --- move the selection to where the author left it ---
$ws.Range("M12").Select()
